# Add 3D sound roll-off fields (area_close / area_far) to the Sound schema
# tables: SOUND_BUNDLE (sheet1) gets the two new columns fully populated
# (header name row, type row, and the three data rows); SOUND_RESOURCE
# (sheet2) only gets the schema header/type cells mirrored (no sample data
# filled in for its two data rows), matching the authored change.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("SOUND_BUNDLE")
$ws2 = $wb.Worksheets.Item("SOUND_RESOURCE")

# ---------------------------------------------------------------------
# SOUND_BUNDLE (sheet1): existing columns are A..F, new columns are G/H
# ---------------------------------------------------------------------

# Row 2 holds the internal field names (index/sound_id/path/loop/volume/
# cooltime) styled like column A's name cell. Clone that formatting onto
# G2:H2, then fill in the two new field names.
$ws1.Range("A2").Copy()
$ws1.Range("G2:H2").PasteSpecial(-4122)
$ws1.Range("G2").Value = "area_close"
$ws1.Range("H2").Value = "area_far"

# Row 3 holds the field types; the existing "cooltime" column (F) is also
# a float, so clone its formatting onto G3:H3 and mark both new fields as
# float too.
$ws1.Range("F3").Copy()
$ws1.Range("G3:H3").PasteSpecial(-4122)
$ws1.Range("G3").Value = "float"
$ws1.Range("H3").Value = "float"

# Data rows: plain numeric roll-off distances for each of the three sound
# bundle entries (no special formatting, matching the other plain data
# columns).
$ws1.Range("G5").Value = 20
$ws1.Range("H5").Value = 50
$ws1.Range("G6").Value = 20
$ws1.Range("H6").Value = 50
$ws1.Range("G7").Value = 20
$ws1.Range("H7").Value = 50

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# SOUND_RESOURCE (sheet2): existing columns are A..G, new columns are H/I
# (schema header only, no sample data rows filled in)
# ---------------------------------------------------------------------

$ws2.Range("A2").Copy()
$ws2.Range("H2:I2").PasteSpecial(-4122)
$ws2.Range("H2").Value = "area_close"
$ws2.Range("I2").Value = "area_far"

$ws2.Range("G3").Copy()
$ws2.Range("H3:I3").PasteSpecial(-4122)
$ws2.Range("H3").Value = "float"
$ws2.Range("I3").Value = "float"

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Switch the active tab back to SOUND_BUNDLE, leaving the cursor on the
# newly added roll-off columns of the last data row.
# ---------------------------------------------------------------------

[void]$ws1.Activate()
[void]$ws1.Range("G6:H6").Select()
